$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cromRange = $ws.Range("C2:C101")
$cromRange.NumberFormat = "@"

$ws.Range("C2").Value = "111101010110010011100010001110"
$ws.Range("D2").Value = 0.9188576210412724
$ws.Range("E2").Value = 0.00425570537100406
$ws.Range("F2").Value = 0.3396782688321202
$ws.Range("C3").Value = "111101010110011001110110011001"
$ws.Range("D3").Value = 0.9189038065527291
$ws.Range("E3").Value = 0.1611554438929474
$ws.Range("F3").Value = 0.8067060084366467
$ws.Range("C4").Value = "111101010110011001110110011100"
$ws.Range("D4").Value = 0.9189038119092944
$ws.Range("E4").Value = 0.8954404232774411
$ws.Range("F4").Value = 0.9165380815963505
$ws.Range("C5").Value = "111101010110011001110110011100"
$ws.Range("D5").Value = 0.9189038119092944
$ws.Range("E5").Value = 0.9188576210412724
$ws.Range("F5").Value = 0.9188958087541856
$ws.Range("C6").Value = "111101010110011001110110011100"
$ws.Range("D6").Value = 0.9189038119092944
$ws.Range("E6").Value = 0.9189038119092944
$ws.Range("F6").Value = 0.9189038119092944
$ws.Range("C7").Value = "111101010110011001110110011100"
$ws.Range("D7").Value = 0.9189038119092944
$ws.Range("E7").Value = 0.9189038119092944
$ws.Range("F7").Value = 0.9189038119092944
$ws.Range("C8").Value = "111101010110011001110110011100"
$ws.Range("D8").Value = 0.9189038119092944
$ws.Range("E8").Value = 0.9189038119092944
$ws.Range("F8").Value = 0.9189038119092944
$ws.Range("C9").Value = "111101010110011001110110011100"
$ws.Range("D9").Value = 0.9189038119092944
$ws.Range("E9").Value = 0.9189038119092944
$ws.Range("F9").Value = 0.9189038119092944
$ws.Range("C10").Value = "111101010110011001110110011100"
$ws.Range("D10").Value = 0.9189038119092944
$ws.Range("E10").Value = 0.9189038119092944
$ws.Range("F10").Value = 0.9189038119092944
$ws.Range("C11").Value = "111101010110011001110110011100"
$ws.Range("D11").Value = 0.9189038119092944
$ws.Range("E11").Value = 0.9189038119092944
$ws.Range("F11").Value = 0.9189038119092944
$ws.Range("C12").Value = "111101010110011001110110011100"
$ws.Range("D12").Value = 0.9189038119092944
$ws.Range("E12").Value = 0.9189038119092944
$ws.Range("F12").Value = 0.9189038119092944
$ws.Range("C13").Value = "111101010110011001110110011100"
$ws.Range("D13").Value = 0.9189038119092944
$ws.Range("E13").Value = 0.9189038119092944
$ws.Range("F13").Value = 0.9189038119092944
$ws.Range("C14").Value = "111101010110011001110110011100"
$ws.Range("D14").Value = 0.9189038119092944
$ws.Range("E14").Value = 0.9189038119092944
$ws.Range("F14").Value = 0.9189038119092944
$ws.Range("C15").Value = "111101010110011001110110011100"
$ws.Range("D15").Value = 0.9189038119092944
$ws.Range("E15").Value = 0.9189038119092944
$ws.Range("F15").Value = 0.9189038119092944
$ws.Range("C16").Value = "111101010110011001110110011100"
$ws.Range("D16").Value = 0.9189038119092944
$ws.Range("E16").Value = 0.9189038119092944
$ws.Range("F16").Value = 0.9189038119092944
$ws.Range("C17").Value = "111101010110011001110110011100"
$ws.Range("D17").Value = 0.9189038119092944
$ws.Range("E17").Value = 0.9189038119092944
$ws.Range("F17").Value = 0.9189038119092944
$ws.Range("C18").Value = "111101010110011001110110011100"
$ws.Range("D18").Value = 0.9189038119092944
$ws.Range("E18").Value = 0.9189038119092944
$ws.Range("F18").Value = 0.9189038119092944
$ws.Range("C19").Value = "111101010110011001110110011100"
$ws.Range("D19").Value = 0.9189038119092944
$ws.Range("E19").Value = 0.9189038119092944
$ws.Range("F19").Value = 0.9189038119092944
$ws.Range("C20").Value = "111101010110011001110110011100"
$ws.Range("D20").Value = 0.9189038119092944
$ws.Range("E20").Value = 0.9189038119092944
$ws.Range("F20").Value = 0.9189038119092944
$ws.Range("C21").Value = "111101010110011001110110011100"
$ws.Range("D21").Value = 0.9189038119092944
$ws.Range("E21").Value = 0.9189038119092944
$ws.Range("F21").Value = 0.9189038119092944
$ws.Range("C22").Value = "111101010110011001110110011100"
$ws.Range("D22").Value = 0.9189038119092944
$ws.Range("E22").Value = 0.9189038119092944
$ws.Range("F22").Value = 0.9189038119092944
$ws.Range("C23").Value = "111101010110011001110110011100"
$ws.Range("D23").Value = 0.9189038119092944
$ws.Range("E23").Value = 0.9189038119092944
$ws.Range("F23").Value = 0.9189038119092944
$ws.Range("C24").Value = "111101010110011001110110011100"
$ws.Range("D24").Value = 0.9189038119092944
$ws.Range("E24").Value = 0.9189038119092944
$ws.Range("F24").Value = 0.9189038119092944
$ws.Range("C25").Value = "111101010110011001110110011100"
$ws.Range("D25").Value = 0.9189038119092944
$ws.Range("E25").Value = 0.9189038119092944
$ws.Range("F25").Value = 0.9189038119092944
$ws.Range("C26").Value = "111101010110011001110110011100"
$ws.Range("D26").Value = 0.9189038119092944
$ws.Range("E26").Value = 0.9189038119092944
$ws.Range("F26").Value = 0.9189038119092944
$ws.Range("C27").Value = "111101010110011001110110011100"
$ws.Range("D27").Value = 0.9189038119092944
$ws.Range("E27").Value = 0.9189038119092944
$ws.Range("F27").Value = 0.9189038119092944
$ws.Range("C28").Value = "111101010110011001110110011100"
$ws.Range("D28").Value = 0.9189038119092944
$ws.Range("E28").Value = 0.9189038119092944
$ws.Range("F28").Value = 0.9189038119092944
$ws.Range("C29").Value = "111101010110011001110110011100"
$ws.Range("D29").Value = 0.9189038119092944
$ws.Range("E29").Value = 0.9189038119092944
$ws.Range("F29").Value = 0.9189038119092944
$ws.Range("C30").Value = "111101010110011001110110011100"
$ws.Range("D30").Value = 0.9189038119092944
$ws.Range("E30").Value = 0.9189038119092944
$ws.Range("F30").Value = 0.9189038119092944
$ws.Range("C31").Value = "111101010110011001110110011100"
$ws.Range("D31").Value = 0.9189038119092944
$ws.Range("E31").Value = 0.9189038119092944
$ws.Range("F31").Value = 0.9189038119092944
$ws.Range("C32").Value = "111101010110011001110110011100"
$ws.Range("D32").Value = 0.9189038119092944
$ws.Range("E32").Value = 0.9189038119092944
$ws.Range("F32").Value = 0.9189038119092944
$ws.Range("C33").Value = "111101010110011001110110011100"
$ws.Range("D33").Value = 0.9189038119092944
$ws.Range("E33").Value = 0.9189038119092944
$ws.Range("F33").Value = 0.9189038119092944
$ws.Range("C34").Value = "111101010110011001110110011100"
$ws.Range("D34").Value = 0.9189038119092944
$ws.Range("E34").Value = 0.9189038119092944
$ws.Range("F34").Value = 0.9189038119092944
$ws.Range("C35").Value = "111101010110011001110110011100"
$ws.Range("D35").Value = 0.9189038119092944
$ws.Range("E35").Value = 0.9189038119092944
$ws.Range("F35").Value = 0.9189038119092944
$ws.Range("C36").Value = "111101010110011001110110011100"
$ws.Range("D36").Value = 0.9189038119092944
$ws.Range("E36").Value = 0.9189038119092944
$ws.Range("F36").Value = 0.9189038119092944
$ws.Range("C37").Value = "111101010110011001110110011100"
$ws.Range("D37").Value = 0.9189038119092944
$ws.Range("E37").Value = 0.9189038119092944
$ws.Range("F37").Value = 0.9189038119092944
$ws.Range("C38").Value = "111101010110011001110110011100"
$ws.Range("D38").Value = 0.9189038119092944
$ws.Range("E38").Value = 0.9189038119092944
$ws.Range("F38").Value = 0.9189038119092944
$ws.Range("C39").Value = "111101010110011001110110011100"
$ws.Range("D39").Value = 0.9189038119092944
$ws.Range("E39").Value = 0.9189038119092944
$ws.Range("F39").Value = 0.9189038119092944
$ws.Range("C40").Value = "111101010110011001110110011100"
$ws.Range("D40").Value = 0.9189038119092944
$ws.Range("E40").Value = 0.9189038119092944
$ws.Range("F40").Value = 0.9189038119092944
$ws.Range("C41").Value = "111101010110011001110110011100"
$ws.Range("D41").Value = 0.9189038119092944
$ws.Range("E41").Value = 0.9189038119092944
$ws.Range("F41").Value = 0.9189038119092944
$ws.Range("C42").Value = "111101010110011001110110011100"
$ws.Range("D42").Value = 0.9189038119092944
$ws.Range("E42").Value = 0.9189038119092944
$ws.Range("F42").Value = 0.9189038119092944
$ws.Range("C43").Value = "111101010110011001110110011100"
$ws.Range("D43").Value = 0.9189038119092944
$ws.Range("E43").Value = 0.9189038119092944
$ws.Range("F43").Value = 0.9189038119092944
$ws.Range("C44").Value = "111101010110011001110110011100"
$ws.Range("D44").Value = 0.9189038119092944
$ws.Range("E44").Value = 0.9189038119092944
$ws.Range("F44").Value = 0.9189038119092944
$ws.Range("C45").Value = "111101010110011001110110011100"
$ws.Range("D45").Value = 0.9189038119092944
$ws.Range("E45").Value = 0.9189038119092944
$ws.Range("F45").Value = 0.9189038119092944
$ws.Range("C46").Value = "111101010110011001110110011100"
$ws.Range("D46").Value = 0.9189038119092944
$ws.Range("E46").Value = 0.9189038119092944
$ws.Range("F46").Value = 0.9189038119092944
$ws.Range("C47").Value = "111101010110011001110110011100"
$ws.Range("D47").Value = 0.9189038119092944
$ws.Range("E47").Value = 0.9189038119092944
$ws.Range("F47").Value = 0.9189038119092944
$ws.Range("C48").Value = "111101010110011001110110011100"
$ws.Range("D48").Value = 0.9189038119092944
$ws.Range("E48").Value = 0.9189038119092944
$ws.Range("F48").Value = 0.9189038119092944
$ws.Range("C49").Value = "111101010110011001110110011100"
$ws.Range("D49").Value = 0.9189038119092944
$ws.Range("E49").Value = 0.9189038119092944
$ws.Range("F49").Value = 0.9189038119092944
$ws.Range("C50").Value = "111101010110011001110110011100"
$ws.Range("D50").Value = 0.9189038119092944
$ws.Range("E50").Value = 0.9189038119092944
$ws.Range("F50").Value = 0.9189038119092944
$ws.Range("C51").Value = "111101010110011001110110011100"
$ws.Range("D51").Value = 0.9189038119092944
$ws.Range("E51").Value = 0.9189038119092944
$ws.Range("F51").Value = 0.9189038119092944
$ws.Range("C52").Value = "111101010110011001110110011100"
$ws.Range("D52").Value = 0.9189038119092944
$ws.Range("E52").Value = 0.9189038119092944
$ws.Range("F52").Value = 0.9189038119092944
$ws.Range("C53").Value = "111101010110011001110110011100"
$ws.Range("D53").Value = 0.9189038119092944
$ws.Range("E53").Value = 0.9189038119092944
$ws.Range("F53").Value = 0.9189038119092944
$ws.Range("C54").Value = "111101010110011001110110011100"
$ws.Range("D54").Value = 0.9189038119092944
$ws.Range("E54").Value = 0.9189038119092944
$ws.Range("F54").Value = 0.9189038119092944
$ws.Range("C55").Value = "111101010110011001110110011100"
$ws.Range("D55").Value = 0.9189038119092944
$ws.Range("E55").Value = 0.9189038119092944
$ws.Range("F55").Value = 0.9189038119092944
$ws.Range("C56").Value = "111101010110011001110110011100"
$ws.Range("D56").Value = 0.9189038119092944
$ws.Range("E56").Value = 0.9189038119092944
$ws.Range("F56").Value = 0.9189038119092944
$ws.Range("C57").Value = "111101010110011001110110011100"
$ws.Range("D57").Value = 0.9189038119092944
$ws.Range("E57").Value = 0.9189038119092944
$ws.Range("F57").Value = 0.9189038119092944
$ws.Range("C58").Value = "111101010110011001110110011100"
$ws.Range("D58").Value = 0.9189038119092944
$ws.Range("E58").Value = 0.9189038119092944
$ws.Range("F58").Value = 0.9189038119092944
$ws.Range("C59").Value = "111101010110011001110110011100"
$ws.Range("D59").Value = 0.9189038119092944
$ws.Range("E59").Value = 0.9189038119092944
$ws.Range("F59").Value = 0.9189038119092944
$ws.Range("C60").Value = "111101010110011001110110011100"
$ws.Range("D60").Value = 0.9189038119092944
$ws.Range("E60").Value = 0.9189038119092944
$ws.Range("F60").Value = 0.9189038119092944
$ws.Range("C61").Value = "111101010110011001110110011100"
$ws.Range("D61").Value = 0.9189038119092944
$ws.Range("E61").Value = 0.9189038119092944
$ws.Range("F61").Value = 0.9189038119092944
$ws.Range("C62").Value = "111101010110011001110110011100"
$ws.Range("D62").Value = 0.9189038119092944
$ws.Range("E62").Value = 0.9189038119092944
$ws.Range("F62").Value = 0.9189038119092944
$ws.Range("C63").Value = "111101010110011001110110011100"
$ws.Range("D63").Value = 0.9189038119092944
$ws.Range("E63").Value = 0.9189038119092944
$ws.Range("F63").Value = 0.9189038119092944
$ws.Range("C64").Value = "111101010110011001110110011100"
$ws.Range("D64").Value = 0.9189038119092944
$ws.Range("E64").Value = 0.9189038119092944
$ws.Range("F64").Value = 0.9189038119092944
$ws.Range("C65").Value = "111101010110011001110110011100"
$ws.Range("D65").Value = 0.9189038119092944
$ws.Range("E65").Value = 0.9189038119092944
$ws.Range("F65").Value = 0.9189038119092944
$ws.Range("C66").Value = "111101010110011001110110011100"
$ws.Range("D66").Value = 0.9189038119092944
$ws.Range("E66").Value = 0.9189038119092944
$ws.Range("F66").Value = 0.9189038119092944
$ws.Range("C67").Value = "111101010110011001110110011100"
$ws.Range("D67").Value = 0.9189038119092944
$ws.Range("E67").Value = 0.9189038119092944
$ws.Range("F67").Value = 0.9189038119092944
$ws.Range("C68").Value = "111101010110011001110110011100"
$ws.Range("D68").Value = 0.9189038119092944
$ws.Range("E68").Value = 0.9189038119092944
$ws.Range("F68").Value = 0.9189038119092944
$ws.Range("C69").Value = "111101010110011001110110011100"
$ws.Range("D69").Value = 0.9189038119092944
$ws.Range("E69").Value = 0.9189038119092944
$ws.Range("F69").Value = 0.9189038119092944
$ws.Range("C70").Value = "111101010110011001110110011100"
$ws.Range("D70").Value = 0.9189038119092944
$ws.Range("E70").Value = 0.9189038119092944
$ws.Range("F70").Value = 0.9189038119092944
$ws.Range("C71").Value = "111101010110011001110110011100"
$ws.Range("D71").Value = 0.9189038119092944
$ws.Range("E71").Value = 0.9189038119092944
$ws.Range("F71").Value = 0.9189038119092944
$ws.Range("C72").Value = "111101010110011001110110011100"
$ws.Range("D72").Value = 0.9189038119092944
$ws.Range("E72").Value = 0.9189038119092944
$ws.Range("F72").Value = 0.9189038119092944
$ws.Range("C73").Value = "111101010110011001110110011100"
$ws.Range("D73").Value = 0.9189038119092944
$ws.Range("E73").Value = 0.9189038119092944
$ws.Range("F73").Value = 0.9189038119092944
$ws.Range("C74").Value = "111101010110011001110110011100"
$ws.Range("D74").Value = 0.9189038119092944
$ws.Range("E74").Value = 0.9189038119092944
$ws.Range("F74").Value = 0.9189038119092944
$ws.Range("C75").Value = "111101010110011001110110011100"
$ws.Range("D75").Value = 0.9189038119092944
$ws.Range("E75").Value = 0.9189038119092944
$ws.Range("F75").Value = 0.9189038119092944
$ws.Range("C76").Value = "111101010110011001110110011100"
$ws.Range("D76").Value = 0.9189038119092944
$ws.Range("E76").Value = 0.9189038119092944
$ws.Range("F76").Value = 0.9189038119092944
$ws.Range("C77").Value = "111101010110011001110110011100"
$ws.Range("D77").Value = 0.9189038119092944
$ws.Range("E77").Value = 0.9189038119092944
$ws.Range("F77").Value = 0.9189038119092944
$ws.Range("C78").Value = "111101010110011001110110011100"
$ws.Range("D78").Value = 0.9189038119092944
$ws.Range("E78").Value = 0.9189038119092944
$ws.Range("F78").Value = 0.9189038119092944
$ws.Range("C79").Value = "111101010110011001110110011100"
$ws.Range("D79").Value = 0.9189038119092944
$ws.Range("E79").Value = 0.9189038119092944
$ws.Range("F79").Value = 0.9189038119092944
$ws.Range("C80").Value = "111101010110011001110110011100"
$ws.Range("D80").Value = 0.9189038119092944
$ws.Range("E80").Value = 0.9189038119092944
$ws.Range("F80").Value = 0.9189038119092944
$ws.Range("C81").Value = "111101010110011001110110011100"
$ws.Range("D81").Value = 0.9189038119092944
$ws.Range("E81").Value = 0.9189038119092944
$ws.Range("F81").Value = 0.9189038119092944
$ws.Range("C82").Value = "111101010110011001110110011100"
$ws.Range("D82").Value = 0.9189038119092944
$ws.Range("E82").Value = 0.9189038119092944
$ws.Range("F82").Value = 0.9189038119092944
$ws.Range("C83").Value = "111101010110011001110110011100"
$ws.Range("D83").Value = 0.9189038119092944
$ws.Range("E83").Value = 0.9189038119092944
$ws.Range("F83").Value = 0.9189038119092944
$ws.Range("C84").Value = "111101010110011001110110011100"
$ws.Range("D84").Value = 0.9189038119092944
$ws.Range("E84").Value = 0.9189038119092944
$ws.Range("F84").Value = 0.9189038119092944
$ws.Range("C85").Value = "111101010110011001110110011100"
$ws.Range("D85").Value = 0.9189038119092944
$ws.Range("E85").Value = 0.9189038119092944
$ws.Range("F85").Value = 0.9189038119092944
$ws.Range("C86").Value = "111101010110011001110110011100"
$ws.Range("D86").Value = 0.9189038119092944
$ws.Range("E86").Value = 0.9189038119092944
$ws.Range("F86").Value = 0.9189038119092944
$ws.Range("C87").Value = "111101010110011001110110011100"
$ws.Range("D87").Value = 0.9189038119092944
$ws.Range("E87").Value = 0.9189038119092944
$ws.Range("F87").Value = 0.9189038119092944
$ws.Range("C88").Value = "111101010110011001110110011100"
$ws.Range("D88").Value = 0.9189038119092944
$ws.Range("E88").Value = 0.9189038119092944
$ws.Range("F88").Value = 0.9189038119092944
$ws.Range("C89").Value = "111101010110011001110110011100"
$ws.Range("D89").Value = 0.9189038119092944
$ws.Range("E89").Value = 0.9189038119092944
$ws.Range("F89").Value = 0.9189038119092944
$ws.Range("C90").Value = "111101010110011001110110011100"
$ws.Range("D90").Value = 0.9189038119092944
$ws.Range("E90").Value = 0.9189038119092944
$ws.Range("F90").Value = 0.9189038119092944
$ws.Range("C91").Value = "111101010110011001110110011100"
$ws.Range("D91").Value = 0.9189038119092944
$ws.Range("E91").Value = 0.9189038119092944
$ws.Range("F91").Value = 0.9189038119092944
$ws.Range("C92").Value = "111101010110011001110110011100"
$ws.Range("D92").Value = 0.9189038119092944
$ws.Range("E92").Value = 0.9189038119092944
$ws.Range("F92").Value = 0.9189038119092944
$ws.Range("C93").Value = "111101010110011001110110011100"
$ws.Range("D93").Value = 0.9189038119092944
$ws.Range("E93").Value = 0.9189038119092944
$ws.Range("F93").Value = 0.9189038119092944
$ws.Range("C94").Value = "111101010110011001110110011100"
$ws.Range("D94").Value = 0.9189038119092944
$ws.Range("E94").Value = 0.9189038119092944
$ws.Range("F94").Value = 0.9189038119092944
$ws.Range("C95").Value = "111101010110011001110110011100"
$ws.Range("D95").Value = 0.9189038119092944
$ws.Range("E95").Value = 0.9189038119092944
$ws.Range("F95").Value = 0.9189038119092944
$ws.Range("C96").Value = "111101010110011001110110011100"
$ws.Range("D96").Value = 0.9189038119092944
$ws.Range("E96").Value = 0.9189038119092944
$ws.Range("F96").Value = 0.9189038119092944
$ws.Range("C97").Value = "111101010110011001110110011100"
$ws.Range("D97").Value = 0.9189038119092944
$ws.Range("E97").Value = 0.9189038119092944
$ws.Range("F97").Value = 0.9189038119092944
$ws.Range("C98").Value = "111101010110011001110110011100"
$ws.Range("D98").Value = 0.9189038119092944
$ws.Range("E98").Value = 0.9189038119092944
$ws.Range("F98").Value = 0.9189038119092944
$ws.Range("C99").Value = "111101010110011001110110011100"
$ws.Range("D99").Value = 0.9189038119092944
$ws.Range("E99").Value = 0.9189038119092944
$ws.Range("F99").Value = 0.9189038119092944
$ws.Range("C100").Value = "111101010110011001110110011100"
$ws.Range("D100").Value = 0.9189038119092944
$ws.Range("E100").Value = 0.9189038119092944
$ws.Range("F100").Value = 0.9189038119092944
$ws.Range("C101").Value = "111101010110011001110110011100"
$ws.Range("D101").Value = 0.9189038119092944
$ws.Range("E101").Value = 0.9189038119092944
$ws.Range("F101").Value = 0.9189038119092944

$cromRange.Style = "Normal"
